$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 1521
    "F3" = 32
    "F4" = 984
    "F5" = 65
    "F6" = 2406
    "F8" = 1480
    "F9" = 69
    "F10" = 167
    "F12" = 418
    "F14" = 8
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cell in $updates.Keys) {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
